$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns: X_Cor -> Lat, Y_Cor -> Long
$ws.Range("B1").Value = "Lat"
$ws.Range("C1").Value = "Long"

# Update selection to reflect new active cell (G2) as recorded in the saved file
$ws.Range("G2").Select() | Out-Null
